$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update transition-probability matrix cells with refreshed simulation results
# (values recomputed after adding more games / optimized simulate-game logic)

# Row 2
$ws.Range("B2").Value = 0.1903114186851211
$ws.Range("C2").Value = 0.5674740484429066
$ws.Range("J2").Value = 0.006920415224913495
$ws.Range("P2").Value = 0.1557093425605536
$ws.Range("S2").Value = 0.07958477508650519

# Row 3
$ws.Range("B3").Value = 0.0119047619047619
$ws.Range("C3").Value = 0.01785714285714286
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("P3").Value = 0.6964285714285714
$ws.Range("S3").Value = 0.2321428571428572

# Row 4
$ws.Range("J4").Value = 0.06976744186046512
$ws.Range("P4").Value = 0.5348837209302325
$ws.Range("S4").Value = 0.3953488372093023

# Row 6
$ws.Range("B6").Value = 0.08290155440414508
$ws.Range("D6").Value = 0.02072538860103627
$ws.Range("F6").Value = 0.08808290155440414
$ws.Range("J6").Value = 0.2124352331606218
$ws.Range("O6").Value = 0.02590673575129534
$ws.Range("Q6").Value = 0.2124352331606218
$ws.Range("R6").Value = 0.07253886010362694
$ws.Range("S6").Value = 0.2849740932642487

# Row 7
$ws.Range("B7").Value = 0.1279069767441861
$ws.Range("D7").Value = 0.01162790697674419
$ws.Range("F7").Value = 0.06395348837209303
$ws.Range("J7").Value = 0.1162790697674419
$ws.Range("O7").Value = 0.02325581395348837
$ws.Range("Q7").Value = 0.1976744186046512
$ws.Range("R7").Value = 0.08139534883720931
$ws.Range("S7").Value = 0.3779069767441861

# Row 8
$ws.Range("B8").Value = 0.09762532981530343
$ws.Range("D8").Value = 0.01846965699208443
$ws.Range("F8").Value = 0.0712401055408971
$ws.Range("J8").Value = 0.129287598944591
$ws.Range("O8").Value = 0.01846965699208443
$ws.Range("Q8").Value = 0.2163588390501319
$ws.Range("R8").Value = 0.0554089709762533
$ws.Range("S8").Value = 0.3931398416886543

# Row 9
$ws.Range("B9").Value = 0.1257861635220126
$ws.Range("D9").Value = 0.01886792452830189
$ws.Range("F9").Value = 0.05031446540880503
$ws.Range("J9").Value = 0.119496855345912
$ws.Range("O9").Value = 0.01257861635220126
$ws.Range("Q9").Value = 0.220125786163522
$ws.Range("R9").Value = 0.06289308176100629
$ws.Range("S9").Value = 0.389937106918239

# Row 10
$ws.Range("B10").Value = 0.1138279932546374
$ws.Range("D10").Value = 0.02613827993254637
$ws.Range("E10").Value = 0.0008431703204047217
$ws.Range("F10").Value = 0.05986509274873524
$ws.Range("J10").Value = 0.1239460370994941
$ws.Range("O10").Value = 0.01602023608768971
$ws.Range("Q10").Value = 0.2377740303541315
$ws.Range("R10").Value = 0.07419898819561552
$ws.Range("S10").Value = 0.3473861720067454

# Row 11
$ws.Range("G11").Value = 0.1477272727272727
$ws.Range("J11").Value = 0.07954545454545454
$ws.Range("K11").Value = 0.1742424242424243
$ws.Range("L11").Value = 0.5909090909090909
$ws.Range("S11").Value = 0.007575757575757576

# Row 12
$ws.Range("G12").Value = 0.7484276729559748
$ws.Range("J12").Value = 0.1635220125786163
$ws.Range("K12").Value = 0.01886792452830189
$ws.Range("L12").Value = 0.01886792452830189
$ws.Range("S12").Value = 0.05031446540880503

# Row 13
$ws.Range("G13").Value = 0.7037037037037037
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("S13").Value = 0.07407407407407407

# Row 15
$ws.Range("F15").Value = 0.0155440414507772
$ws.Range("H15").Value = 0.1398963730569948
$ws.Range("I15").Value = 0.07253886010362694
$ws.Range("J15").Value = 0.3367875647668394
$ws.Range("K15").Value = 0.06735751295336788
$ws.Range("M15").Value = 0.02590673575129534
$ws.Range("O15").Value = 0.04145077720207254
$ws.Range("S15").Value = 0.3005181347150259

# Row 16
$ws.Range("F16").Value = 0.02857142857142857
$ws.Range("H16").Value = 0.1657142857142857
$ws.Range("I16").Value = 0.05714285714285714
$ws.Range("J16").Value = 0.4742857142857143
$ws.Range("K16").Value = 0.09142857142857143
$ws.Range("M16").Value = 0.02285714285714286
$ws.Range("N16").Value = 0.005714285714285714
$ws.Range("O16").Value = 0.04571428571428571
$ws.Range("S16").Value = 0.1085714285714286

# Row 17
$ws.Range("F17").Value = 0.02536997885835095
$ws.Range("H17").Value = 0.1585623678646934
$ws.Range("I17").Value = 0.08456659619450317
$ws.Range("J17").Value = 0.4418604651162791
$ws.Range("K17").Value = 0.09090909090909091
$ws.Range("M17").Value = 0.01691331923890063
$ws.Range("O17").Value = 0.06553911205073996
$ws.Range("S17").Value = 0.1162790697674419

# Row 18
$ws.Range("F18").Value = 0.006944444444444444
$ws.Range("H18").Value = 0.1805555555555556
$ws.Range("I18").Value = 0.08333333333333333
$ws.Range("J18").Value = 0.4652777777777778
$ws.Range("K18").Value = 0.08333333333333333
$ws.Range("O18").Value = 0.0763888888888889
$ws.Range("S18").Value = 0.1041666666666667

# Row 19
$ws.Range("F19").Value = 0.01555352241537054
$ws.Range("H19").Value = 0.202195791399817
$ws.Range("I19").Value = 0.07502287282708142
$ws.Range("J19").Value = 0.3989021043000915
$ws.Range("K19").Value = 0.1143641354071363
$ws.Range("M19").Value = 0.01097895699908509
$ws.Range("O19").Value = 0.06587374199451052
$ws.Range("S19").Value = 0.1171088746569076
